$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same date serial (45188 = 2023-09-19) for
# every data row (2..431). The update bumps it by one day to 45189
# (2023-09-20) across the whole column.
$ws.Range("C2:C431").Value = 45189
